$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append 4 new lesson rows (rows 56-59) ---------------------------------
$ws.Range("A56").Value = 45051
$ws.Range("B56").Value = "Your First Solo"
$ws.Range("C56").Value = "Drag"

$ws.Range("A57").Value = 45051
$ws.Range("B57").Value = "Your First Solo"
$ws.Range("C57").Value = "A closer look : Reducing Drag"

$ws.Range("A58").Value = 45051
$ws.Range("B58").Value = "Your First Solo"
$ws.Range("C58").Value = "Thrust,Stability, and Center of Gravity"

$ws.Range("A59").Value = 45051
$ws.Range("B59").Value = "Your First Solo"
$ws.Range("C59").Value = "Flight service Weather Briefings"

# Match the date formatting used by the rest of column A (reuse the existing
# style instead of creating a brand-new number format xf).
$ws.Range("A55").Copy()
$ws.Range("A56:A59").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# --- Update the sheet selection --------------------------------------------
$ws.Cells.Select()
